$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 855.40814
$ws.Range("I15").Value = 855.40814
$ws.Range("K15").Value = 2566.22442
$ws.Range("M15").Value = -2397.22442

$ws.Range("H39").Value = 756.6667
$ws.Range("I39").Value = 175.55556
$ws.Range("J39").Value = 2500
$ws.Range("K39").Value = 526.66668
$ws.Range("L39").Value = 7500
$ws.Range("M39").Value = -230.66668
$ws.Range("N39").Value = -8092

$ws.Range("H53").Value = 518.36365
$ws.Range("I53").Value = 258.2143
$ws.Range("K53").Value = 258.2143
$ws.Range("M53").Value = 378.7857

$ws.Range("H64").Value = 5002.5293
$ws.Range("J64").Value = 5855.778
$ws.Range("L64").Value = 5855.778
$ws.Range("N64").Value = -6351.778

$ws.Range("H67").Value = 5002.5293
$ws.Range("J67").Value = 5855.778
$ws.Range("L67").Value = 5855.778
$ws.Range("N67").Value = -7571.778

$ws.Range("H74").Value = 5900
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -7872

$ws.Range("H77").Value = 5900
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -39360

$ws.Range("H135").Value = 7643.3335
$ws.Range("I135").Value = 4472.75
$ws.Range("J135").Value = 10179.8
$ws.Range("K135").Value = 40254.75
$ws.Range("L135").Value = 91618.2
$ws.Range("M135").Value = -37719.75
$ws.Range("N135").Value = -96688.2

$ws.Range("H137").Value = 3801.8
$ws.Range("I137").Value = 2236.75
$ws.Range("K137").Value = 6710.25
$ws.Range("M137").Value = -4160.25

$ws.Range("H138").Value = 6936.132
$ws.Range("J138").Value = 9602.362999999999
$ws.Range("L138").Value = 28807.089
$ws.Range("N138").Value = -39087.089

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 879.75
$ws.Range("I5").Value = 879.75
$ws.Range("K5").Value = 879.75
$ws.Range("M5").Value = -767.75

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H32").Value = 1909.9012
$ws.Range("I32").Value = 1032.1466
$ws.Range("K32").Value = 1032.1466
$ws.Range("M32").Value = -745.1466

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 879.75
$ws.Range("I4").Value = 879.75
$ws.Range("K4").Value = 879.75
$ws.Range("M4").Value = -764.75

$ws.Range("H86").Value = 8501.723
$ws.Range("J86").Value = 2999.6667
$ws.Range("L86").Value = 2999.6667
$ws.Range("N86").Value = -5245.6667

$ws.Range("H89").Value = 8501.723
$ws.Range("J89").Value = 2999.6667
$ws.Range("L89").Value = 14998.3335
$ws.Range("N89").Value = -26230.3335

$ws.Range("H134").Value = 2535.3076
$ws.Range("I134").Value = 2400.95
$ws.Range("J134").Value = 2983.1667
$ws.Range("K134").Value = 7202.849999999999
$ws.Range("L134").Value = 8949.500100000001
$ws.Range("M134").Value = -4667.849999999999
$ws.Range("N134").Value = -14019.5001

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 22584
$ws.Range("I2").Value = 38.4
$ws.Range("J2").Value = 60160
$ws.Range("K2").Value = 38.4
$ws.Range("L2").Value = 60160
$ws.Range("M2").Value = 74.59999999999999
$ws.Range("N2").Value = -60386

$ws.Range("H7").Value = 187.09375
$ws.Range("I7").Value = 102.85714
$ws.Range("K7").Value = 102.85714
$ws.Range("M7").Value = 10.14286

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H31").Value = 2097.8076
$ws.Range("J31").Value = 2293.2173
$ws.Range("L31").Value = 2293.2173
$ws.Range("N31").Value = -2883.2173

$ws.Range("H34").Value = 2097.8076
$ws.Range("J34").Value = 2293.2173
$ws.Range("L34").Value = 2293.2173
$ws.Range("N34").Value = -2697.2173

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H99").Value = 10771.292
$ws.Range("J99").Value = 12443.3125
$ws.Range("L99").Value = 12443.3125
$ws.Range("N99").Value = -15439.3125

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H126").Value = 10771.292
$ws.Range("J126").Value = 12443.3125
$ws.Range("L126").Value = 37329.9375
$ws.Range("N126").Value = -42269.9375

$ws.Range("H132").Value = 50686.34
$ws.Range("I132").Value = 55820.055
$ws.Range("K132").Value = 167460.165
$ws.Range("M132").Value = -164930.165

$ws.Range("H134").Value = 4604.9346
$ws.Range("I134").Value = 4707.8096
$ws.Range("J134").Value = 3524.75
$ws.Range("K134").Value = 14123.4288
$ws.Range("L134").Value = 10574.25
$ws.Range("M134").Value = -11588.4288
$ws.Range("N134").Value = -15644.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1387900.6
$ws.Range("J32").Value = 2020241
$ws.Range("L32").Value = 6060723
$ws.Range("N32").Value = -6061289

$ws.Range("H55").Value = 4455868.5
$ws.Range("J55").Value = 6067514.5
$ws.Range("L55").Value = 18202543.5
$ws.Range("N55").Value = -18202897.5

$ws.Range("H117").Value = 1384.2858
$ws.Range("I117").Value = 627.55554
$ws.Range("J117").Value = 2746.4
$ws.Range("K117").Value = 1882.66662
$ws.Range("L117").Value = 8239.200000000001
$ws.Range("M117").Value = 1559.33338
$ws.Range("N117").Value = -15123.2

$ws.Range("H132").Value = 5162
$ws.Range("I132").Value = 3680.25
$ws.Range("J132").Value = 6643.75
$ws.Range("K132").Value = 33122.25
$ws.Range("L132").Value = 59793.75
$ws.Range("M132").Value = -30592.25
$ws.Range("N132").Value = -64853.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1428901.9
$ws.Range("I2").Value = 5555766
$ws.Range("J2").Value = 372.03845
$ws.Range("K2").Value = 5555766
$ws.Range("L2").Value = 372.03845
$ws.Range("M2").Value = -5555653
$ws.Range("N2").Value = -598.03845

$ws.Range("H122").Value = 3850.3333
$ws.Range("I122").Value = 3399.875
$ws.Range("J122").Value = 4751.25
$ws.Range("K122").Value = 10199.625
$ws.Range("L122").Value = 14253.75
$ws.Range("M122").Value = -7749.625
$ws.Range("N122").Value = -19153.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3000
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -3340

$ws.Range("H22").Value = 2984.8696
$ws.Range("I22").Value = 2631.5557
$ws.Range("J22").Value = 4256.8
$ws.Range("K22").Value = 2631.5557
$ws.Range("L22").Value = 4256.8
$ws.Range("M22").Value = -2336.5557
$ws.Range("N22").Value = -4846.8

$ws.Range("H27").Value = 2984.8696
$ws.Range("I27").Value = 2631.5557
$ws.Range("J27").Value = 4256.8
$ws.Range("K27").Value = 2631.5557
$ws.Range("L27").Value = 4256.8
$ws.Range("M27").Value = -2524.5557
$ws.Range("N27").Value = -4470.8

$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5344

$ws.Range("H40").Value = 4226.909
$ws.Range("I40").Value = 4267.4614
$ws.Range("J40").Value = 4076.2856
$ws.Range("K40").Value = 4267.4614
$ws.Range("L40").Value = 4076.2856
$ws.Range("M40").Value = -4131.4614
$ws.Range("N40").Value = -4348.2856

$ws.Range("H46").Value = 1918.3429
$ws.Range("I46").Value = 979.2593000000001
$ws.Range("K46").Value = 979.2593000000001
$ws.Range("M46").Value = -791.2593000000001

$ws.Range("H122").Value = 4066.4167
$ws.Range("I122").Value = 3696.6667
$ws.Range("J122").Value = 4682.6665
$ws.Range("K122").Value = 11090.0001
$ws.Range("L122").Value = 14047.9995
$ws.Range("M122").Value = -8640.000100000001
$ws.Range("N122").Value = -18947.9995

$ws.Range("H136").Value = 1542212.1
$ws.Range("I136").Value = 2003539.2
$ws.Range("K136").Value = 6010617.6
$ws.Range("M136").Value = -6008067.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 28495
$ws.Range("J33").Value = 28495
$ws.Range("L33").Value = 28495
$ws.Range("N33").Value = -28995

$ws.Range("H36").Value = 28495
$ws.Range("J36").Value = 28495
$ws.Range("L36").Value = 28495
$ws.Range("N36").Value = -28995

$ws.Range("H37").Value = 20833
$ws.Range("J37").Value = 20833
$ws.Range("L37").Value = 20833
$ws.Range("N37").Value = -21239

$ws.Range("H75").Value = 53332.332
$ws.Range("J75").Value = 54999
$ws.Range("L75").Value = 54999
$ws.Range("N75").Value = -56871

$ws.Range("H78").Value = 53332.332
$ws.Range("J78").Value = 54999
$ws.Range("L78").Value = 164997
$ws.Range("N78").Value = -174357

$ws.Range("H136").Value = 2546.1282
$ws.Range("J136").Value = 2640.875
$ws.Range("L136").Value = 7922.625
$ws.Range("N136").Value = -13022.625
